$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowValues = @{
    "A" = "2021年"
    "C" = 1024
    "E" = 879
    "F" = 1211
    "G" = 3114
    "H" = 21
    "J" = 147695
    "K" = 4621
    "L" = 12441
    "M" = 482
    "N" = 377
    "O" = 22364
    "P" = 1929
    "R" = 654
    "S" = 3087
    "T" = 302
    "U" = 2416
    "V" = 1754
    "W" = 33
    "X" = 11626
    "Y" = 1429
    "AA" = 13455
    "AB" = 625
    "AC" = 249
    "AD" = 3458
    "AE" = 2940
    "AF" = 2
    "AH" = 5520
    "AI" = 12469
    "AJ" = 1382
    "AK" = 1795
    "AL" = 1070
    "AM" = 960
    "AN" = 8518
    "AO" = 2440
    "AP" = 1391
    "AQ" = 3483
    "AR" = 1
    "AS" = 1362
    "AT" = 9990
    "AU" = 138
    "AV" = 1024
    "AW" = 5564
    "AX" = 95
}

$emptyCols = @("B", "D", "I", "Q", "Z", "AG")

foreach ($col in $rowValues.Keys) {
    $cellRef = "$($col)12"
    $ws.Range($cellRef).Value = $rowValues[$col]
}

foreach ($col in $emptyCols) {
    $cellRef = "$($col)12"
    $ws.Range($cellRef).Value = "'"
    $ws.Range($cellRef).Style = "Normal"
}

$ws.Range("A11").Copy()
$ws.Range("A12").PasteSpecial(-4122)
